# Added Q-tail (column N, "Q_union") pixel-width values for each glyph row.
# (Comma-tail / column O "comma_union" is left for a follow-up edit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> Q_union value, mirroring the lookup table's existing "q_union" (M) column data
$qUnionValues = @{
    2=59; 3=118; 4=59; 5=119; 6=59; 7=120; 8=59; 9=118; 10=121; 11=119;
    12=118; 13=122; 14=59; 15=59; 16=59; 17=59; 18=59; 19=59; 20=59; 21=121;
    22=59; 23=59; 24=59; 25=59; 26=59; 27=59; 28=123; 29=124; 30=123; 31=124;
    32=125; 33=125; 34=123; 35=126; 36=125; 37=127; 38=126; 39=118; 40=128; 41=126;
    42=123; 43=124; 44=123; 45=124; 46=123; 47=125; 48=126; 49=126; 50=129; 51=126;
    52=126; 53=125; 54=123; 55=121; 56=123; 57=123; 58=119; 59=125; 60=127; 61=125;
    62=123; 63=123; 64=59; 65=59; 66=121; 67=121; 68=59; 69=130; 70=131; 71=132;
    72=133; 73=59; 74=122; 75=127
}

foreach ($row in ($qUnionValues.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 14).Value = $qUnionValues[$row]
}

# The comma row (row 65, ",") got its Q_union cell underlined as a visual
# reminder that the comma-tail column is the next thing to fill in.
$ws.Cells.Item(65, 14).Font.Underline = $true

# Restore the view: scroll down a little further and land the selection on G55
# inside the frozen-pane (right) half of the sheet, same as the saved workbook.
$ws.Range("G55").Select() | Out-Null
